$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.908.83"
$ws.Range("D3").Value = "1.815.09"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'309.29"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4653"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.3663"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'0.8689"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "1.860.88"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "'0.07090"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "'6.513"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "'0.000008715"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "26.938.49"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'5.301"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "'10.65"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "2.068.75"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").Value = "'1.892"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").Value = "'150.99"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").Value = "'18.31"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "'2.140"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "'5.257"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "'115.39"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "'0.08918"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'0.7570"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").Value = "'1.157"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "'4.487"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "'2.910"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "'1.084"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").Value = "'0.05276"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'0.01947"
$ws.Range("D40").Value = "'2.980"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").Value = "'7.244"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'0.5301"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'2.293"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").Value = "'0.1652"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "'8.434"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "'0.4879"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'103.30"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").Value = "'1.661"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "'0.06291"
$ws.Range("E51").Value = "  +0.13%  "
